$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing covid_deaths (column C) values
$ws.Cells.Item(741,3).Value = 4
$ws.Cells.Item(813,3).Value = 5
$ws.Cells.Item(912,3).Value = 25
$ws.Cells.Item(1036,3).Value = 39
$ws.Cells.Item(1044,3).Value = 32
$ws.Cells.Item(1050,3).Value = 38
$ws.Cells.Item(1053,3).Value = 13
$ws.Cells.Item(1085,3).Value = 51
$ws.Cells.Item(1090,3).Value = 35
$ws.Cells.Item(1116,3).Value = 25
$ws.Cells.Item(1117,3).Value = 36
$ws.Cells.Item(1122,3).Value = 19
$ws.Cells.Item(1128,3).Value = 20
$ws.Cells.Item(1129,3).Value = 28
$ws.Cells.Item(1131,3).Value = 6
$ws.Cells.Item(1132,3).Value = 18
$ws.Cells.Item(1133,3).Value = 23
$ws.Cells.Item(1136,3).Value = 4
$ws.Cells.Item(1138,3).Value = 37
$ws.Cells.Item(1139,3).Value = 4
$ws.Cells.Item(1140,3).Value = 13
$ws.Cells.Item(1141,3).Value = 14
$ws.Cells.Item(1142,3).Value = 30
$ws.Cells.Item(1144,3).Value = 2
$ws.Cells.Item(1145,3).Value = 9
$ws.Cells.Item(1146,3).Value = 23
$ws.Cells.Item(1147,3).Value = 29
$ws.Cells.Item(1148,3).Value = 2
$ws.Cells.Item(1149,3).Value = 6
$ws.Cells.Item(1150,3).Value = 7
$ws.Cells.Item(1151,3).Value = 16
$ws.Cells.Item(1152,3).Value = 28

# Append new rows for date 44194 (2020-12-29)
$ws.Cells.Item(1153,1).Value = 44194
$ws.Cells.Item(1153,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1153,2).Value = "40-49"
$ws.Cells.Item(1153,3).Value = 1
$ws.Cells.Item(1154,1).Value = 44194
$ws.Cells.Item(1154,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1154,2).Value = "50-59"
$ws.Cells.Item(1154,3).Value = 3
$ws.Cells.Item(1155,1).Value = 44194
$ws.Cells.Item(1155,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1155,2).Value = "60-69"
$ws.Cells.Item(1155,3).Value = 3
$ws.Cells.Item(1156,1).Value = 44194
$ws.Cells.Item(1156,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1156,2).Value = "70-79"
$ws.Cells.Item(1156,3).Value = 10
$ws.Cells.Item(1157,1).Value = 44194
$ws.Cells.Item(1157,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1157,2).Value = "80+"
$ws.Cells.Item(1157,3).Value = 14
